# ---------------------------------------------------------------------------
# Adds two new columns (Admin Fee, Total Amount) after "Funding Amount",
# turns each employee's dependent info into its own "Dependent:" sub-row,
# and recomputes the Funding/Admin Fee/Total figures + the trailing Total row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that LOOKS like a number but must stay plain text
# (matches the source data, which stores "3.00" / "0.20" / ... as text).
# Forcing text format first keeps Excel's auto-detection from turning the
# string into a real number, then ClearFormats() drops the temporary
# NumberFormat so the cell ends up back on the sheet's default style.
function Set-TextValue($rangeRef, $text) {
    $rng = $ws.Range($rangeRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# 1) Insert two new blank columns at I:J. Excel shifts the old
#    I/J/K ("Dependent Name"/"Dependent Tier"/"Dependent Relationship")
#    to K/L/M automatically, carrying their values and header styling.
$ws.Columns("I:J").Insert()

# 2) New column headers.
$ws.Range("I1").Value = "Admin Fee"
$ws.Range("J1").Value = "Total Amount"

# 3) Row 2 (Jane Johnson's primary row): new funding/fee/total, and the
#    dependent columns on this row are cleared (dependent moves to row 3).
Set-TextValue "H2" "3.00"
Set-TextValue "I2" "0.20"
Set-TextValue "J2" "3.20"
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""

# 4) Row 3 becomes the "Dependent:" sub-row for Jane Johnson: employee
#    identity columns (GL Code/Location/Title/Plan/Tier) are cleared,
#    and the funding/fee/total reflect the dependent's own amounts.
#    K3/L3/M3 already hold "Casey Smith"/"Tier1"/"Spouse" post-insert.
$ws.Range("A3").Value = "Dependent:"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""
Set-TextValue "H3" "39.00"
Set-TextValue "I3" "9.10"
Set-TextValue "J3" "48.10"

# 5) Rows 4-7 (Alice Smith, Alice Brown, Mike Doe, Alex Williams): same
#    Funding/Admin Fee/Total pattern, no dependents.
Set-TextValue "H4" "3.00"
Set-TextValue "I4" "0.20"
Set-TextValue "J4" "3.20"

Set-TextValue "H5" "3.00"
Set-TextValue "I5" "0.20"
Set-TextValue "J5" "3.20"

Set-TextValue "H6" "3.00"
Set-TextValue "I6" "0.20"
Set-TextValue "J6" "3.20"

Set-TextValue "H7" "3.00"
Set-TextValue "I7" "0.20"
Set-TextValue "J7" "3.20"

# 6) Row 8 totals.
$ws.Range("A8").Value = "Total:"
Set-TextValue "H8" "54.00"
Set-TextValue "I8" "10.10"
Set-TextValue "J8" "64.10"
